$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")

$ws1.Cells.Item(2,4).Value = 3290.05
$ws1.Cells.Item(2,5).Value = 102.24
$ws1.Cells.Item(3,4).Value = 2965
$ws1.Cells.Item(3,5).Value = 765
$ws1.Cells.Item(4,4).Value = 2700
$ws1.Cells.Item(4,5).Value = 690
$ws1.Cells.Item(5,1).Value = "SETAO CI"
$ws1.Cells.Item(5,4).Value = 2630
$ws1.Cells.Item(5,5).Value = 680
$ws1.Cells.Item(6,1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(6,4).Value = 2548.33
$ws1.Cells.Item(6,5).Value = 632.62
$ws1.Cells.Item(7,5).Value = 600
$ws1.Cells.Item(8,4).Value = 2350
$ws1.Cells.Item(9,4).Value = 1494.48
$ws1.Cells.Item(9,5).Value = 380.09
$ws1.Cells.Item(10,4).Value = 1433.14
$ws1.Cells.Item(10,5).Value = 346.35
$ws1.Cells.Item(11,4).Value = 1319.77
$ws1.Cells.Item(11,5).Value = 330.94
$ws1.Cells.Item(12,1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(12,3).Value = 4
$ws1.Cells.Item(12,4).Value = 562.45
$ws1.Cells.Item(12,5).Value = 138.46
$ws1.Cells.Item(13,1).Value = "BRVM - INDUSTRIE"
$ws1.Cells.Item(13,3).Value = 2
$ws1.Cells.Item(13,4).Value = 538.85
$ws1.Cells.Item(13,5).Value = 269.52
$ws1.Cells.Item(14,1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(14,3).Value = 4
$ws1.Cells.Item(14,4).Value = 527.8
$ws1.Cells.Item(14,5).Value = 132
$ws1.Cells.Item(15,1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(15,4).Value = 497.83
$ws1.Cells.Item(15,5).Value = 124.36
$ws1.Cells.Item(16,1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(16,4).Value = 489.27
$ws1.Cells.Item(16,5).Value = 122.22
$ws1.Cells.Item(17,1).Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Cells.Item(17,3).Value = 2
$ws1.Cells.Item(17,4).Value = 444.05
$ws1.Cells.Item(17,5).Value = 222.19
$ws1.Cells.Item(18,1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(18,4).Value = 439.07
$ws1.Cells.Item(18,5).Value = 109.74
$ws1.Cells.Item(19,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(19,4).Value = 438.08
$ws1.Cells.Item(19,5).Value = 112.8
$ws1.Cells.Item(20,1).Value = "BRVM-PRINCIPAL"
$ws1.Cells.Item(20,3).Value = 2
$ws1.Cells.Item(20,4).Value = 385.88
$ws1.Cells.Item(20,5).Value = 193.48
$ws1.Cells.Item(21,4).Value = 376.31
$ws1.Cells.Item(21,5).Value = 94.43000000000001
$ws1.Cells.Item(23,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(23,3).Value = 1
$ws1.Cells.Item(23,4).Value = 14.32
$ws1.Cells.Item(23,5).Value = 7.24
$ws1.Cells.Item(24,1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(24,2).Value = 2
$ws1.Cells.Item(24,4).Value = 11.01
$ws1.Cells.Item(24,5).Value = 6.98
$ws1.Cells.Item(25,1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(25,2).Value = 2
$ws1.Cells.Item(25,3).Value = 1
$ws1.Cells.Item(25,4).Value = 8.08
$ws1.Cells.Item(25,5).Value = 4.17
$ws1.Cells.Item(25,7).Value = "👀 À surveiller"
$ws1.Cells.Item(26,1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(26,4).Value = 7.25
$ws1.Cells.Item(26,5).Value = 7.25
$ws1.Cells.Item(27,1).Value = "CIE CI (CIEC)"
$ws1.Cells.Item(27,2).Value = 1
$ws1.Cells.Item(27,3).Value = 0
$ws1.Cells.Item(27,4).Value = 4.74
$ws1.Cells.Item(27,5).Value = 4.74
$ws1.Cells.Item(27,7).Value = "➖ Neutre"
$ws1.Cells.Item(28,1).Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Cells.Item(28,4).Value = 3.48
$ws1.Cells.Item(28,5).Value = 3.48
$ws1.Cells.Item(29,3).Value = 0
$ws1.Cells.Item(29,4).Value = 3.16
$ws1.Cells.Item(29,7).Value = "➖ Neutre"
$ws1.Cells.Item(30,1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(30,3).Value = 0
$ws1.Cells.Item(30,4).Value = 2.86
$ws1.Cells.Item(30,5).Value = 2.86
$ws1.Cells.Item(30,7).Value = "➖ Neutre"
$ws1.Cells.Item(31,1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(31,3).Value = 0
$ws1.Cells.Item(31,4).Value = 2.74
$ws1.Cells.Item(31,5).Value = 2.74
$ws1.Cells.Item(31,7).Value = "➖ Neutre"
$ws1.Cells.Item(32,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(32,4).Value = 0.37
$ws1.Cells.Item(32,5).Value = 6.25
$ws1.Cells.Item(33,1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(33,2).Value = 1
$ws1.Cells.Item(33,3).Value = 1
$ws1.Cells.Item(33,4).Value = 0.25
$ws1.Cells.Item(33,5).Value = 5.09
$ws1.Cells.Item(33,7).Value = "👀 À surveiller"
$ws1.Cells.Item(34,1).Value = "TOTAL"
$ws1.Cells.Item(34,2).Value = 0
$ws1.Cells.Item(34,4).Value = 0
$ws1.Cells.Item(34,5).Value = 0
$ws1.Cells.Item(34,7).Value = "➖ Neutre"
$ws1.Cells.Item(35,1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(35,4).Value = -1.71
$ws1.Cells.Item(35,5).Value = -1.71
$ws1.Cells.Item(36,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(36,4).Value = -2
$ws1.Cells.Item(36,5).Value = -2
$ws1.Cells.Item(37,1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(37,4).Value = -2.13
$ws1.Cells.Item(37,5).Value = -2.13
$ws1.Cells.Item(41,1).Value = "SITAB CI (STBC)"
$ws1.Cells.Item(41,2).Value = 0
$ws1.Cells.Item(41,3).Value = 1
$ws1.Cells.Item(41,4).Value = -3.04
$ws1.Cells.Item(41,5).Value = -3.04
$ws1.Cells.Item(41,7).Value = "➖ Neutre"
$ws1.Cells.Item(42,1).Value = "SUCRIVOIRE (SCRC)"
$ws1.Cells.Item(42,2).Value = 0
$ws1.Cells.Item(42,3).Value = 1
$ws1.Cells.Item(42,4).Value = -3.28
$ws1.Cells.Item(42,5).Value = -3.28
$ws1.Cells.Item(42,7).Value = "➖ Neutre"
$ws1.Cells.Item(43,1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Cells.Item(43,3).Value = 1
$ws1.Cells.Item(43,4).Value = -3.47
$ws1.Cells.Item(43,5).Value = -3.47
$ws1.Cells.Item(44,1).Value = "BICI CI (BICC)"
$ws1.Cells.Item(44,2).Value = 0
$ws1.Cells.Item(44,3).Value = 2
$ws1.Cells.Item(44,4).Value = -3.98
$ws1.Cells.Item(44,5).Value = -2.5
$ws1.Cells.Item(44,6).Value = "🟡 Observer"
$ws1.Cells.Item(44,7).Value = "➖ Neutre"
$ws1.Cells.Item(45,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(45,2).Value = 1
$ws1.Cells.Item(45,3).Value = 2
$ws1.Cells.Item(45,4).Value = -4.9
$ws1.Cells.Item(45,5).Value = -4.76
$ws1.Cells.Item(45,6).Value = "🟡 Observer"
$ws1.Cells.Item(45,7).Value = "👀 À surveiller"
$ws1.Cells.Item(46,1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(46,2).Value = 0
$ws1.Cells.Item(46,3).Value = 3
$ws1.Cells.Item(46,4).Value = -8.58
$ws1.Cells.Item(46,5).Value = -3.51
$ws1.Cells.Item(46,6).Value = "🔴 Vente"
$ws1.Cells.Item(46,7).Value = "⚠️ Risque de décrochage"

$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Cells.Item(2,2).Value = 7527675.26
$ws2.Cells.Item(3,2).Value = 500400.89
$ws2.Cells.Item(4,2).Value = 360515.65
$ws2.Cells.Item(5,1).Value = "SETAO CI"
$ws2.Cells.Item(5,2).Value = 328134.34
$ws2.Cells.Item(6,1).Value = "BRVM - AUTRES SECTEURS"
$ws2.Cells.Item(6,2).Value = 295055.39
$ws2.Cells.Item(8,2).Value = 223298.02
$ws2.Cells.Item(9,2).Value = 50204.49
$ws2.Cells.Item(10,2).Value = 43925.46
$ws2.Cells.Item(11,2).Value = 34069.54
